$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

# The "#system" sheet is hidden; make it visible so it can be manipulated,
# then restore its hidden state at the end.
$sys.Visible = $true

# ---------------------------------------------------------------------------
# 1) Remove the lone "text" category (column Y, a single-member list) and
#    close the gap by shifting columns Z:AE one column to the left (into
#    Y:AD). A whole-column delete is safe here since it only touches Y:AE.
# ---------------------------------------------------------------------------
$sys.Range("Y:Y").Delete()

# ---------------------------------------------------------------------------
# 2) The master category list in column A ("target") also listed "text" at
#    A25; remove that single cell and shift A26:A31 up into A25:A30 --
#    scoped strictly to column A (cell-by-cell, since Range.Delete on a
#    partial column shifts the whole row in this host).
# ---------------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $val = $sys.Cells.Item($r + 1, 1).Value2
    $sys.Cells.Item($r, 1).Value = $val
}
$sys.Cells.Item(31, 1).ClearContents()

# ---------------------------------------------------------------------------
# 3) Add the new json function `storeKeys(json,jsonpath,var)` to column M,
#    inserted alphabetically before `storeValue` (M16) -- shift M16:M17 down
#    into M17:M18, scoped strictly to column M (cell-by-cell, same reason).
# ---------------------------------------------------------------------------
for ($r = 17; $r -ge 16; $r--) {
    $val = $sys.Cells.Item($r, 13).Value2
    $sys.Cells.Item($r + 1, 13).Value = $val
}
$sys.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 4) Update the named ranges to reflect the new extents / shifted columns.
# ---------------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"

$sys.Visible = $false
